# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Re-orders / refreshes the worker mora data block (rows 16-21) on sheet
# "Hoja1": three workers (MARILYN CANTILLO AVILA, MARTHA IRINA TEJEDOR
# PANZA, ROSA ALEJANDRA CARRASQUILLA RONCALLO), each now listed once for
# period 1712 and once for period 1801 (grouped by period instead of by
# worker), and refreshes ROSA's "Salario Basico" value for both periods.
# Also nudges the company logo image left to match the new column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row data: row, TipoDoc, NDoc, Nombre, PeriodoMora, ValorMora, SalarioBasico
$rows = @(
    @(16, "CC", "1143327174", "MARILYN CANTILLO AVILA",                "1712", 29509, 737717),
    @(17, "CC", "1148434693", "MARTHA IRINA TEJEDOR PANZA",             "1712", 29509, 737717),
    @(18, "CC", "1128054808", "ROSA ALEJANDRA CARRASQUILLA RONCALLO",   "1712", 40920, 1023000),
    @(19, "CC", "1143327174", "MARILYN CANTILLO AVILA",                "1801", 29509, 737717),
    @(20, "CC", "1148434693", "MARTHA IRINA TEJEDOR PANZA",             "1801", 29509, 737717),
    @(21, "CC", "1128054808", "ROSA ALEJANDRA CARRASQUILLA RONCALLO",   "1801", 40920, 1023000)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value2 = $r[1]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($rowNum, 3).Value2 = $r[2]   # C - N Doc Trabajador
    $ws.Cells.Item($rowNum, 4).Value2 = $r[3]   # D - Nombre Trabajador
    $ws.Cells.Item($rowNum, 5).Value2 = $r[4]   # E - Periodo Mora
    $ws.Cells.Item($rowNum, 6).Value2 = $r[5]   # F - Valor Mora
    $ws.Cells.Item($rowNum, 7).Value2 = $r[6]   # G - Salario Basico
}

# --- Columns B:J were re-autofit by Excel after the data refresh; apply the
#     resulting best-fit widths (characters).
$ws.Columns.Item(2).ColumnWidth = 16.90625
$ws.Columns.Item(3).ColumnWidth = 10.81640625
$ws.Columns.Item(4).ColumnWidth = 38.26953125
$ws.Columns.Item(5).ColumnWidth = 12.7265625
$ws.Columns.Item(6).ColumnWidth = 9.453125
$ws.Columns.Item(7).ColumnWidth = 13.453125
$ws.Columns.Item(8).ColumnWidth = 17.90625
$ws.Columns.Item(9).ColumnWidth = 16.81640625
$ws.Columns.Item(10).ColumnWidth = 14.1796875

# --- The logo picture shifted 13.5pt to the left (column B got narrower).
$logo = $ws.Shapes.Item(1)
$logo.Left = $logo.Left - 13.5
